$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 37662.332
$ws.Range("J47").Value = 41194.8
$ws.Range("L47").Value = 41194.8
$ws.Range("N47").Value = -43138.8
$ws.Range("H112").Value = 2698.225
$ws.Range("I112").Value = 1846.625
$ws.Range("J112").Value = 2911.125
$ws.Range("K112").Value = 5539.875
$ws.Range("L112").Value = 8733.375
$ws.Range("M112").Value = -4431.875
$ws.Range("N112").Value = -10949.375
$ws.Range("H125").Value = 1869.8334
$ws.Range("J125").Value = 2208.6667
$ws.Range("L125").Value = 19878.0003
$ws.Range("N125").Value = -24798.0003

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 18890.223
$ws.Range("I45").Value = 20726.5
$ws.Range("K45").Value = 20726.5
$ws.Range("M45").Value = -20349.5
$ws.Range("H61").Value = 2071
$ws.Range("J61").Value = 1999.5
$ws.Range("L61").Value = 1999.5
$ws.Range("N61").Value = -2423.5
$ws.Range("H136").Value = 2071
$ws.Range("J136").Value = 1999.5
$ws.Range("L136").Value = 5998.5
$ws.Range("N136").Value = -11098.5
$ws.Range("H139").Value = 79232.44500000001
$ws.Range("J139").Value = 79232.44500000001
$ws.Range("L139").Value = 79232.44500000001
$ws.Range("N139").Value = -89512.44500000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 200000
$ws.Range("J42").Value = 200000
$ws.Range("L42").Value = 200000
$ws.Range("N42").Value = -200656
$ws.Range("H99").Value = 4483.5654
$ws.Range("I99").Value = 4246
$ws.Range("J99").Value = 5338.8
$ws.Range("K99").Value = 4246
$ws.Range("L99").Value = 5338.8
$ws.Range("M99").Value = -2748
$ws.Range("N99").Value = -8334.799999999999
$ws.Range("H105").Value = 1878.8572
$ws.Range("I105").Value = 1881.5227
$ws.Range("J105").Value = 1855.4
$ws.Range("K105").Value = 1881.5227
$ws.Range("L105").Value = 1855.4
$ws.Range("M105").Value = -134.5227
$ws.Range("N105").Value = -5349.4

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 6508.55
$ws.Range("J4").Value = 6588.1055
$ws.Range("L4").Value = 6588.1055
$ws.Range("N4").Value = -6812.1055
$ws.Range("H58").Value = 10581.235
$ws.Range("I58").Value = 11455.5
$ws.Range("J58").Value = 8978.416999999999
$ws.Range("K58").Value = 11455.5
$ws.Range("L58").Value = 8978.416999999999
$ws.Range("M58").Value = -11252.5
$ws.Range("N58").Value = -9384.416999999999
$ws.Range("H136").Value = 10581.235
$ws.Range("I136").Value = 11455.5
$ws.Range("J136").Value = 8978.416999999999
$ws.Range("K136").Value = 34366.5
$ws.Range("L136").Value = 26935.251
$ws.Range("M136").Value = -31816.5
$ws.Range("N136").Value = -32035.251

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 513.29034
$ws.Range("I5").Value = 525.4
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 1576.2
$ws.Range("L5").Value = 450
$ws.Range("M5").Value = -1464.2
$ws.Range("N5").Value = -674
$ws.Range("H8").Value = 2859.6667
$ws.Range("I8").Value = 2859.6667
$ws.Range("K8").Value = 8579.000100000001
$ws.Range("M8").Value = -8440.000100000001
$ws.Range("H135").Value = 513.29034
$ws.Range("I135").Value = 525.4
$ws.Range("J135").Value = 150
$ws.Range("K135").Value = 4728.599999999999
$ws.Range("L135").Value = 1350
$ws.Range("M135").Value = -2193.599999999999
$ws.Range("N135").Value = -6420
$ws.Range("H137").Value = 3884.7144
$ws.Range("I137").Value = 4598.6
$ws.Range("J137").Value = 2100
$ws.Range("K137").Value = 13795.8
$ws.Range("L137").Value = 6300
$ws.Range("M137").Value = -8695.800000000001
$ws.Range("N137").Value = -16500

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 899
$ws.Range("I4").Value = 899
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 899
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("M4").Value = -787
$ws.Range("H70").Value = 16455.79
$ws.Range("J70").Value = 6724.5713
$ws.Range("L70").Value = 6724.5713
$ws.Range("N70").Value = -7264.5713
$ws.Range("H73").Value = 16455.79
$ws.Range("J73").Value = 6724.5713
$ws.Range("L73").Value = 6724.5713
$ws.Range("N73").Value = -8596.5713
$ws.Range("H113").Value = 3505.2
$ws.Range("I113").Value = 4173
$ws.Range("J113").Value = 3145.6155
$ws.Range("K113").Value = 4173
$ws.Range("L113").Value = 3145.6155
$ws.Range("M113").Value = -2003
$ws.Range("N113").Value = -7485.6155
$ws.Range("H122").Value = 3519.6155
$ws.Range("I122").Value = 3579.5715
$ws.Range("K122").Value = 10738.7145
$ws.Range("M122").Value = -8288.7145
$ws.Range("H126").Value = 4685.524
$ws.Range("I126").Value = 4109.933
$ws.Range("K126").Value = 12329.799
$ws.Range("M126").Value = -9859.798999999999
$ws.Range("H132").Value = 8571.723
$ws.Range("I132").Value = 7679.2
$ws.Range("K132").Value = 23037.6
$ws.Range("M132").Value = -20507.6

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8885.842000000001
$ws.Range("I7").Value = 8379.385
$ws.Range("J7").Value = 9983.166999999999
$ws.Range("K7").Value = 8379.385
$ws.Range("L7").Value = 9983.166999999999
$ws.Range("M7").Value = -8267.385
$ws.Range("N7").Value = -10207.167
$ws.Range("H22").Value = 1585.2142
$ws.Range("I22").Value = 1847.25
$ws.Range("J22").Value = 1235.8334
$ws.Range("K22").Value = 1847.25
$ws.Range("L22").Value = 1235.8334
$ws.Range("M22").Value = -1552.25
$ws.Range("N22").Value = -1825.8334
$ws.Range("H27").Value = 1585.2142
$ws.Range("I27").Value = 1847.25
$ws.Range("J27").Value = 1235.8334
$ws.Range("K27").Value = 1847.25
$ws.Range("L27").Value = 1235.8334
$ws.Range("M27").Value = -1740.25
$ws.Range("N27").Value = -1449.8334
$ws.Range("H40").Value = 2820.5806
$ws.Range("I40").Value = 2682.4614
$ws.Range("K40").Value = 2682.4614
$ws.Range("M40").Value = -2546.4614
$ws.Range("H61").Value = 11971.12
$ws.Range("I61").Value = 12377.739
$ws.Range("J61").Value = 7295
$ws.Range("K61").Value = 12377.739
$ws.Range("L61").Value = 7295
$ws.Range("M61").Value = -12175.739
$ws.Range("N61").Value = -7699
$ws.Range("H68").Value = 2949.6191
$ws.Range("I68").Value = 2960.3572
$ws.Range("J68").Value = 2928.1428
$ws.Range("K68").Value = 2960.3572
$ws.Range("L68").Value = 2928.1428
$ws.Range("M68").Value = -2211.3572
$ws.Range("N68").Value = -4426.1428
$ws.Range("H71").Value = 2949.6191
$ws.Range("I71").Value = 2960.3572
$ws.Range("J71").Value = 2928.1428
$ws.Range("K71").Value = 14801.786
$ws.Range("L71").Value = 14640.714
$ws.Range("M71").Value = -11057.786
$ws.Range("N71").Value = -22128.714
$ws.Range("H113").Value = 11971.12
$ws.Range("I113").Value = 12377.739
$ws.Range("J113").Value = 7295
$ws.Range("K113").Value = 12377.739
$ws.Range("L113").Value = 7295
$ws.Range("M113").Value = -10207.739
$ws.Range("N113").Value = -11635
$ws.Range("H126").Value = 8885.842000000001
$ws.Range("I126").Value = 8379.385
$ws.Range("J126").Value = 9983.166999999999
$ws.Range("K126").Value = 25138.155
$ws.Range("L126").Value = 29949.501
$ws.Range("M126").Value = -22668.155
$ws.Range("N126").Value = -34889.501
$ws.Range("H136").Value = 7206.136
$ws.Range("I136").Value = 2808.5625
$ws.Range("K136").Value = 8425.6875
$ws.Range("M136").Value = -5875.6875
$ws.Range("H141").Value = 161622
$ws.Range("J141").Value = 161622
$ws.Range("L141").Value = 161622
$ws.Range("N141").Value = -171982

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 50000
$ws.Range("J47").Value = 50000
$ws.Range("L47").Value = 50000
$ws.Range("N47").Value = -51144
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("L48").Value = 0
$ws.Range("H126").Value = 2797.8696
$ws.Range("I126").Value = 2874.45
$ws.Range("K126").Value = 8623.349999999999
$ws.Range("M126").Value = -6153.349999999999
